# Add the new "count-the-smiley-faces" column (G) and fix up a handful of
# "tests-results" (F) values that were previously recorded incorrectly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in G1, matching the style of the other header cells (F1).
$ws.Range("G1").Value = "count-the-smiley-faces"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats

# New column G values for each data row (2-27).
$gValues = @{
    2  = "No"
    3  = "Good"
    4  = "No"
    5  = "No"
    6  = "No"
    7  = "Good"
    8  = "No"
    9  = "Good"
    10 = "Good"
    11 = "Good"
    12 = "No"
    13 = "Good"
    14 = "Good"
    15 = "No"
    16 = "No"
    17 = "No"
    18 = "No"
    19 = "Good"
    20 = "Good"
    21 = "Good"
    22 = "Good"
    23 = "No"
    24 = "No"
    25 = "Good"
    26 = "No"
    27 = "Good"
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}

# A handful of existing "tests-results" (F) cells were corrected to "Good".
$fFixes = @(6, 9, 18, 24)
foreach ($row in $fFixes) {
    $ws.Cells.Item($row, 6).Value = "Good"
}
